$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-CellValue($rowIndex, $colIndex, $value) {
    $row = $tbl.Rows.Item($rowIndex)
    $cell = $row.Cells.Item($colIndex)
    $r = $cell.Range
    $r.Text = $value
    $r.Font.Bold = $true
    $r.Font.Size = 12
    $r.Font.SizeBi = 12
}

# RETENTION table "Ratio" row
Set-CellValue 24 2 "0.3333"

# Answer Recall Lenient (ARL)
Set-CellValue 44 2 "0.0833"

# Answer Recall Strict (ARS)
Set-CellValue 45 2 "0"

# Answer Recall Average (ARA)
Set-CellValue 46 2 "0.0416"
